$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "58.029.51"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "2.456.53"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'511.01"
$ws.Range("E5").Value = "  -2.50%  "
$ws.Range("D6").Value = "'133.81"
$ws.Range("E6").Value = "  +3.02%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'0.558"
$ws.Range("E8").Value = "  -1.44%  "
$ws.Range("D9").Value = "2.456.91"
$ws.Range("D10").Value = "'0.0983"
$ws.Range("E10").Value = "  +0.79%  "
$ws.Range("E11").Value = "  -0.49%  "
$ws.Range("D13").Value = "'4.63"
$ws.Range("E13").Value = "  -6.84%  "
$ws.Range("D14").Value = "2.890.51"
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("D15").Value = "57.899.77"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").Value = "'21.94"
$ws.Range("E16").Value = "  +1.74%  "
$ws.Range("D17").Value = "'0.0000135"
$ws.Range("E17").Value = "  +1.49%  "
$ws.Range("D18").Value = "2.449.87"
$ws.Range("E18").Value = "  -0.08%  "
$ws.Range("D19").Value = "'10.35"
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("D21").Value = "'315.23"
$ws.Range("E21").Value = "  +0.68%  "
$ws.Range("E22").Value = "  +5.02%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "'5.73"
$ws.Range("E24").Value = "  -1.83%  "
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.37%  "
$ws.Range("E27").Value = "  -0.73%  "
$ws.Range("E28").Value = "  -5.81%  "
$ws.Range("D29").Value = "'7.58"
$ws.Range("E29").Value = "  +4.70%  "
$ws.Range("D30").Value = "'172.75"
$ws.Range("E30").Value = "  -1.29%  "
$ws.Range("D31").Value = "0.0₃0734"
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("E32").Value = "  +0.35%  "
$ws.Range("D33").Value = "'6.18"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").Value = "'1.13"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("E37").Value = "  +1.54%  "
$ws.Range("E38").Value = "  +5.40%  "
$ws.Range("E39").Value = "  +2.36%  "
$ws.Range("E40").Value = "  +1.31%  "
$ws.Range("E41").Value = "  +1.47%  "
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("D43").Value = "'137.00"
$ws.Range("E43").Value = "  +9.35%  "
$ws.Range("E44").Value = "  +0.99%  "
$ws.Range("D45").Value = "'4.92"
$ws.Range("E45").Value = "  +2.73%  "
$ws.Range("E46").Value = "  -1.17%  "
$ws.Range("D47").Value = "'256.79"
$ws.Range("E47").Value = "  -0.52%  "
$ws.Range("D48").Value = "'0.0920"
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("D49").Value = "'0.0494"
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("E50").Value = "  +1.91%  "
$ws.Range("D51").Value = "'17.19"
$ws.Range("E51").Value = "  +0.93%  "
